$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

# 1. "document as sets" -> "documents as sets"
$d.Content.Find.Execute("represent document as sets", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "represent documents as sets", $wdReplaceAll)

# 2. remove duplicated whitespace before "the goal" (scoped so the earlier
#    run's <w:tab/> is left untouched)
$d.Content.Find.Execute("  the goal", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, " the goal", $wdReplaceAll)

# 3. add missing period at end of whitespace-shingling sentence
$d.Content.Find.Execute("any-long whitespace sequence as 1 character", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "any-long whitespace sequence as 1 character.", $wdReplaceAll)

# 4. "which hold the characteristic vectors" -> "which holds the characteristic vectors"
$d.Content.Find.Execute("which hold the characteristic vectors", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "which holds the characteristic vectors", $wdReplaceAll)

# remove duplicate "holds" in "The characteristic vector holds compares"
$d.Content.Find.Execute("characteristic vector holds compares", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "characteristic vector compares", $wdReplaceAll)

# 5. "takes too much resources" -> "demands too many resources" (scoped to avoid other "takes"/"much")
$rngTakes = $d.Content
$rngTakes.Find.Execute("but that also takes too much resources. Instead,")
$subTakes = $d.Range($rngTakes.Start, $rngTakes.End)
$subTakes.Find.Execute("takes", $true, $false, $false, $false, $false,
                        $true, $wdFindContinue, $false, "demands", $wdReplaceAll)
$subTakes2 = $d.Range($rngTakes.Start, $rngTakes.End)
$subTakes2.Find.Execute("much", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "many", $wdReplaceAll)

# 6. remove "The shingling logic can be found in the Shingler class. " sentence
$d.Content.Find.Execute("Java. The shingling logic can be found in the Shingler class. First", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "Java. First", $wdReplaceAll)

# 7. "(1 entry)" -> "(1 medical entry)"
$d.Content.Find.Execute("document (1 entry) based", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "document (1 medical entry) based", $wdReplaceAll)

# 8. interface/swing framework paragraph edits
$d.Content.Find.Execute("my application I used the swing framework", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "my application, I used the java swing framework", $wdReplaceAll)
$d.Content.Find.Execute("shingling size, before which", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "shingling size before which", $wdReplaceAll)

# 9. rewrite ending about extra features
$d.Content.Find.Execute("the highest or the few highest estimated similarities or exact similarities.", $true, $false, $false, $false, $false,
                         $true, $wdFindContinue, $false, "the highest exact or estimated similarities or generating similar pairs based on estimates.", $wdReplaceAll)
